$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.286989212036133
$ws.Range("B1").Value = 3.542252540588379
$ws.Range("C1").Value = 1.75562310218811
$ws.Range("D1").Value = 1.196251153945923
$ws.Range("E1").Value = 1.256491541862488
